# feat: add 2022-Q4 data
#
# Plan:
#  1. The existing "2022-Q3" sheet (position 2) holds the quarter's fund-holdings
#     table. We duplicate it first so the ORIGINAL values survive unmodified as the
#     new, final "2022-Q3" sheet (position 3).
#  2. The original sheet object (still at position 2) is then turned into the new
#     "2022-Q4" sheet: renamed, and its data + header/index-cell styling updated to
#     the new quarter's numbers (style copied over from the "总计" sheet's own
#     header formatting, matching the target workbook).
#  3. The summary sheet "总计" gets its existing "2022-Q3" label swapped to
#     "2022-Q4" and a brand-new appended row for "2022-Q3" (duplicating the prior
#     row's counts, matching the target).
#
# Note: some numbers in the source data (e.g. "0.21", "86.58") are stored as
# plain TEXT, not numeric cells. Assigning a numeric-looking string straight to
# `.Value` makes Excel auto-convert it to a real number, so instead we stage the
# text in a helper cell pre-formatted as Text ("@"), copy it, and
# PasteSpecial-values it into the destination — that carries over only the
# value (as text), not any number formatting, onto the target cell.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# --- Step 1: duplicate the current "2022-Q3" sheet so its original data is kept
# as-is in a sheet placed right after it (this will end up being the final,
# unchanged "2022-Q3" sheet). ---
$q3Sheet.Copy($null, $q3Sheet)
$q3Copy = $wb.ActiveSheet
$q3Copy.Name = "2022-Q3-NEW-TMP"

# --- Step 2: turn the original sheet into the new "2022-Q4" sheet. ---
$q3Sheet.Name = "2022-Q4"
$q4Sheet = $q3Sheet

# Helper cell used to push numeric-looking strings in as literal TEXT (see
# note above) instead of letting Excel coerce them into numbers.
$helper = $totalSheet.Range("Z100")
$helper.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

# Update the fund rows with the 2022-Q4 figures.
Set-TextValue $q4Sheet.Range("D2") "0.21"
Set-TextValue $q4Sheet.Range("E2") "86.58"
Set-TextValue $q4Sheet.Range("F2") "5.37"
Set-TextValue $q4Sheet.Range("G2") "0.0113"
$q4Sheet.Range("H2").Value = 5

Set-TextValue $q4Sheet.Range("D3") "0.16"
Set-TextValue $q4Sheet.Range("E3") "86.58"
Set-TextValue $q4Sheet.Range("F3") "5.37"
Set-TextValue $q4Sheet.Range("G3") "0.0086"
$q4Sheet.Range("H3").Value = 5

$helper.Clear()

# Re-style the header row + index column to match the target workbook (copied
# from the "总计" sheet's own header / first index cell formatting).
$totalSheet.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A3").PasteSpecial(-4122)

# --- Step 3: rename the duplicated sheet back to "2022-Q3" (its data is already
# identical to the original, untouched). ---
$q3Copy.Name = "2022-Q3"

# --- Step 4: update the "总计" summary sheet. ---
$totalSheet.Range("B2").Value = "2022-Q4"

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.02

# Restore the original active sheet ("总计" was selected before any edits).
$totalSheet.Activate()
